$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 103
$ws.Range("H103").Value = 45455164
$ws.Range("I103").Value = 639
$ws.Range("J103").Value = 125000584
$ws.Range("K103").Value = 1917
$ws.Range("L103").Value = 375001752
$ws.Range("M103").Value = -1331
$ws.Range("N103").Value = -375002924
# Row 107
$ws.Range("H107").Value = 1371.92
$ws.Range("I107").Value = 1427.7778
$ws.Range("J107").Value = 1228.2858
$ws.Range("K107").Value = 1427.7778
$ws.Range("L107").Value = 1228.2858
$ws.Range("M107").Value = 492.2221999999999
$ws.Range("N107").Value = -5068.2858
# Row 132
$ws.Range("H132").Value = 1382.3235
$ws.Range("I132").Value = 1349.9688
$ws.Range("K132").Value = 4049.9064
$ws.Range("M132").Value = -1519.9064
# Row 138
$ws.Range("H138").Value = 3756.3447
$ws.Range("I138").Value = 1372.35
$ws.Range("J138").Value = 9054.111000000001
$ws.Range("K138").Value = 4117.049999999999
$ws.Range("L138").Value = 27162.333
$ws.Range("M138").Value = 1022.950000000001
$ws.Range("N138").Value = -37442.333

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 1496.5769
$ws.Range("I74").Value = 1496.5769
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1496.5769
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -622.5769
$ws.Range("N74").Value = ""
# Row 77
$ws.Range("H77").Value = 1496.5769
$ws.Range("I77").Value = 1496.5769
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 7482.8845
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -3114.8845
$ws.Range("N77").Value = ""

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1554.4546
$ws.Range("I99").Value = 909.9
$ws.Range("J99").Value = 8000
$ws.Range("K99").Value = 909.9
$ws.Range("L99").Value = 8000
$ws.Range("M99").Value = 588.1
$ws.Range("N99").Value = -10996
# Row 105
$ws.Range("H105").Value = 1581.8182
$ws.Range("I105").Value = 1198.0769
$ws.Range("K105").Value = 1198.0769
$ws.Range("M105").Value = 548.9231
# Row 107
$ws.Range("H107").Value = 201958.6
$ws.Range("I107").Value = 1630.6666
$ws.Range("K107").Value = 1630.6666
$ws.Range("M107").Value = 289.3334

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1574.0454
$ws.Range("I16").Value = 1385.6428
$ws.Range("J16").Value = 1903.75
$ws.Range("K16").Value = 1385.6428
$ws.Range("L16").Value = 1903.75
$ws.Range("M16").Value = -1098.6428
$ws.Range("N16").Value = -2477.75
# Row 23
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = ""
# Row 27
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = ""
# Row 58
$ws.Range("H58").Value = 1682.8148
$ws.Range("I58").Value = 1571.5
$ws.Range("J58").Value = 2000.8572
$ws.Range("K58").Value = 1571.5
$ws.Range("L58").Value = 2000.8572
$ws.Range("M58").Value = -1368.5
$ws.Range("N58").Value = -2406.8572
# Row 105
$ws.Range("H105").Value = 2067
$ws.Range("I105").Value = 1534
$ws.Range("J105").Value = 2600
$ws.Range("K105").Value = 1534
$ws.Range("L105").Value = 2600
$ws.Range("M105").Value = 213
$ws.Range("N105").Value = -6094
# Row 113
$ws.Range("H113").Value = 1574.0454
$ws.Range("I113").Value = 1385.6428
$ws.Range("J113").Value = 1903.75
$ws.Range("K113").Value = 1385.6428
$ws.Range("L113").Value = 1903.75
$ws.Range("M113").Value = 784.3571999999999
$ws.Range("N113").Value = -6243.75
# Row 132
$ws.Range("H132").Value = 2446.353
$ws.Range("I132").Value = 2245.3333
$ws.Range("J132").Value = 3221.7144
$ws.Range("K132").Value = 6735.999899999999
$ws.Range("L132").Value = 9665.143199999999
$ws.Range("M132").Value = -4205.999899999999
$ws.Range("N132").Value = -14725.1432
# Row 134
$ws.Range("H134").Value = 2093.4783
$ws.Range("I134").Value = 2007.1428
$ws.Range("K134").Value = 6021.428400000001
$ws.Range("M134").Value = -3486.428400000001
# Row 136
$ws.Range("H136").Value = 1682.8148
$ws.Range("I136").Value = 1571.5
$ws.Range("J136").Value = 2000.8572
$ws.Range("K136").Value = 4714.5
$ws.Range("L136").Value = 6002.571599999999
$ws.Range("M136").Value = -2164.5
$ws.Range("N136").Value = -11102.5716

$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 2870.9678
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2870.9678
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 8612.903399999999
$ws.Range("M80").Value = ""
$ws.Range("N80").Value = -10484.9034
# Row 83
$ws.Range("H83").Value = 2870.9678
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2870.9678
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 25838.7102
$ws.Range("M83").Value = ""
$ws.Range("N83").Value = -35198.7102
# Row 131
$ws.Range("H131").Value = 767.0599999999999
$ws.Range("I131").Value = 398.33334
$ws.Range("J131").Value = 817.3409
$ws.Range("K131").Value = 1195.00002
$ws.Range("L131").Value = 2452.0227
$ws.Range("M131").Value = 3844.99998
$ws.Range("N131").Value = -12532.0227

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 83334990
$ws.Range("I113").Value = 1372.8334
$ws.Range("J113").Value = 166668620
$ws.Range("K113").Value = 1372.8334
$ws.Range("L113").Value = 166668620
$ws.Range("M113").Value = 797.1666
$ws.Range("N113").Value = -166672960

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 2390
$ws.Range("I61").Value = 1254
$ws.Range("J61").Value = 4175.143
$ws.Range("K61").Value = 1254
$ws.Range("L61").Value = 4175.143
$ws.Range("M61").Value = -1052
$ws.Range("N61").Value = -4579.143
# Row 113
$ws.Range("H113").Value = 2390
$ws.Range("I113").Value = 1254
$ws.Range("J113").Value = 4175.143
$ws.Range("K113").Value = 1254
$ws.Range("L113").Value = 4175.143
$ws.Range("M113").Value = 916
$ws.Range("N113").Value = -8515.143
# Row 136
$ws.Range("H136").Value = 11498457
$ws.Range("I136").Value = 16669620
$ws.Range("J136").Value = 6982.6665
$ws.Range("K136").Value = 50008860
$ws.Range("L136").Value = 20947.9995
$ws.Range("M136").Value = -50006310
$ws.Range("N136").Value = -26047.9995

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 1158
$ws.Range("I107").Value = 798.3333
$ws.Range("J107").Value = 1589.6
$ws.Range("K107").Value = 2394.9999
$ws.Range("L107").Value = 4768.799999999999
$ws.Range("M107").Value = -474.9998999999998
$ws.Range("N107").Value = -8608.799999999999
# Row 113
$ws.Range("H113").Value = 25000426
$ws.Range("I113").Value = 421.7
$ws.Range("J113").Value = 100000430
$ws.Range("K113").Value = 1265.1
$ws.Range("L113").Value = 300001290
$ws.Range("M113").Value = 904.9000000000001
$ws.Range("N113").Value = -300005630
# Row 122
$ws.Range("H122").Value = 1629.8096
$ws.Range("I122").Value = 1668.0769
$ws.Range("J122").Value = 1567.625
$ws.Range("K122").Value = 5004.2307
$ws.Range("L122").Value = 4702.875
$ws.Range("M122").Value = -2554.2307
$ws.Range("N122").Value = -9602.875
# Row 132
$ws.Range("H132").Value = 24196362
$ws.Range("I132").Value = 32609822
$ws.Range("J132").Value = 7670.5625
$ws.Range("K132").Value = 97829466
$ws.Range("L132").Value = 23011.6875
$ws.Range("M132").Value = -97826936
$ws.Range("N132").Value = -28071.6875
# Row 136
$ws.Range("H136").Value = 1463.125
$ws.Range("I136").Value = 1250.8334
$ws.Range("K136").Value = 3752.5002
$ws.Range("M136").Value = -1202.5002
